$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with new recipe names
$ws.Range("B1").Value = "Recipe Chicken in creme freice"
$ws.Range("C1").Value = "Recipe Chilli con carne"
$ws.Range("D1").Value = "Recipe Kylling og søtpotet"
$ws.Range("E1").Value = "Recipe Meeeeeat"
$ws.Range("F1").Value = "Recipe Porrige"

# Update data values for rows 2-7 (columns B-F)
$ws.Range("B2").Value = 0.07090242000000001
$ws.Range("C2").Value = 0.09367762000000002
$ws.Range("D2").Value = 0.06318120000000001
$ws.Range("E2").Value = 0.3644722
$ws.Range("F2").Value = 0.02128775

$ws.Range("B3").Value = 0.036352752
$ws.Range("C3").Value = 0.070164752
$ws.Range("D3").Value = 0.03232159999999999
$ws.Range("E3").Value = 0.3354869999999999
$ws.Range("F3").Value = 0.01837625

$ws.Range("B4").Value = 844.6276000000001
$ws.Range("C4").Value = 1321.8376
$ws.Range("D4").Value = 462.52
$ws.Range("E4").Value = 2044.39
$ws.Range("F4").Value = 646.3399999999999

$ws.Range("B5").Value = 29224.90480000001
$ws.Range("C5").Value = 37519.82480000001
$ws.Range("D5").Value = 15544.88
$ws.Range("E5").Value = 59785.87000000001
$ws.Range("F5").Value = 23344.01

$ws.Range("B6").Value = 7.9744924
$ws.Range("C6").Value = 17.5627824
$ws.Range("D6").Value = 6.32596
$ws.Range("E6").Value = 110.73899
$ws.Range("F6").Value = 5.075285

$ws.Range("B7").Value = 12.1432608
$ws.Range("C7").Value = 23.6602708
$ws.Range("D7").Value = 7.74464
$ws.Range("E7").Value = 354.65746
$ws.Range("F7").Value = 10.51516

# Remove column G entirely (data + used-range shrinks back to A1:F7)
$ws.Range("G1:G7").Delete()
